$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.918.29"
$ws.Range("E2").Value = "  +3.05%  "

$ws.Range("D3").Value = "1.570.00"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").Value = "'0.995"
$ws.Range("E4").Value = "  -1.34%  "

$ws.Range("D5").Value = "'211.38"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "'0.492"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -1.24%  "

$ws.Range("D8").Value = "'23.13"
$ws.Range("E8").Value = "  +5.40%  "

$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("D10").Value = "'0.0595"
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").Value = "'0.0880"
$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").Value = "1.796.80"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").Value = "1.562.36"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").Value = "'3.74"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").Value = "'0.519"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "27.914.30"
$ws.Range("E16").Value = "  +3.22%  "

$ws.Range("D17").Value = "'63.29"
$ws.Range("E17").Value = "  +1.99%  "

$ws.Range("D18").Value = "'227.75"
$ws.Range("E18").Value = "  +5.75%  "

$ws.Range("D19").Value = "0.0₃0703"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "'7.42"
$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").Value = "'4.09"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("D23").Value = "'9.25"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").Value = "'151.56"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "'15.14"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").Value = "'6.56"
$ws.Range("E27").Value = "  -0.69%  "

$ws.Range("D28").Value = "'0.106"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").Value = "'3.21"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("D33").Value = "'3.12"
$ws.Range("E33").Value = "  -1.88%  "

$ws.Range("D34").Value = "1.408.94"
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("E35").Value = "  -1.50%  "

$ws.Range("E36").Value = "  -3.52%  "

$ws.Range("E37").Value = "  -1.92%  "

$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").Value = "'0.538"
$ws.Range("E39").Value = "  +1.38%  "

$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("D41").Value = "'0.803"
$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("D43").Value = "'5.60"
$ws.Range("E43").Value = "  -3.69%  "

$ws.Range("D44").Value = "'0.971"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("E45").Value = "  +4.16%  "

$ws.Range("D46").Value = "'63.41"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("D47").Value = "1.707.64"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "'86.51"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0524"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("E51").Value = "  -1.57%  "
